$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.551.50"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "2.491.62"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'494.00"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'153.21"
$ws.Range("E6").Value = "  +7.28%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "2.502.64"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("E10").Value = "  +4.54%  "
$ws.Range("D11").Value = "'0.0988"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D14").Value = "2.926.90"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "56.709.12"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").Value = "'21.42"
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "2.507.44"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("D20").Value = "'10.33"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").Value = "'321.79"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "'0.997"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "'5.89"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").Value = "'59.00"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("E27").Value = "  -3.00%  "
$ws.Range("D28").Value = "2.615.15"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +1.67%  "
$ws.Range("D30").Value = "0.0₃0813"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "'151.75"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").Value = "'18.36"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("D34").Value = "'1.52"
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").Value = "'0.869"
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("E39").Value = "  +4.19%  "
$ws.Range("D40").Value = "'33.94"
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("E41").Value = "  +2.12%  "
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").Value = "'0.617"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").Value = "'4.93"
$ws.Range("E45").Value = "  +5.20%  "
$ws.Range("D46").Value = "'269.50"
$ws.Range("E46").Value = "  +5.91%  "
$ws.Range("D47").Value = "'0.0927"
$ws.Range("E47").Value = "  +1.89%  "
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").Value = "'10.22"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").Value = "'17.85"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").Value = "1.895.30"
$ws.Range("E51").Value = "  -5.61%  "
